## "1. added post status"
## Adds two new status/config sheets - "ccpd_post_st" and "uc_running_cfg" -
## right after "config" and before "heater_cfg", then makes "ccpd_post_st"
## the active tab.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. Create the two new sheets in an order that reproduces the target
#        sheetId allocation (ccpd_post_st=12, uc_running_cfg=11): create a
#        throwaway filler first (takes id 10), then uc_running_cfg (11),
#        then ccpd_post_st (12), then drop the filler. All three land right
#        before "heater_cfg" (position 2), giving a final tab order of:
#        config, ccpd_post_st, uc_running_cfg, heater_cfg, ...
$filler = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$filler.Name = "zz_filler_zz"

$ucRunning = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$ucRunning.Name = "uc_running_cfg"

$ccpdPost = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$ccpdPost.Name = "ccpd_post_st"

$wb.Worksheets.Item("zz_filler_zz").Delete() | Out-Null

$header = @("filename", "marker", "common_id", "field", "c_datatype", "bytes", "msg")

# --- 2. Populate "uc_running_cfg" first (its brand-new strings occupy the
#        shared-string table before ccpd_post_st's do). -------------------
$ws2 = $wb.Worksheets.Item("uc_running_cfg")

for ($i = 0; $i -lt $header.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $header[$i]
}
$ws2.Range("A1:G1").Font.Bold = $true

$ws2.Cells.Item(2, 1).Value = "utilities/uc_helper.h"
$ws2.Cells.Item(2, 2).Value = "UC_RUNNING_CFG_STRUCT"

$ws2.Cells.Item(2, 4).Value = "uc_type"
$ws2.Cells.Item(3, 4).Value = "uc_state_id"
$ws2.Cells.Item(4, 4).Value = "uc_state_action"

$ws2.Cells.Item(2, 5).Value = "uint8_t"
$ws2.Cells.Item(3, 5).Value = "uint16_t"
$ws2.Cells.Item(4, 5).Value = "int8_t"

$ws2.Cells.Item(2, 3).Value = 100
$ws2.Cells.Item(2, 6).Value = 1
$ws2.Cells.Item(3, 6).Value = 2
$ws2.Cells.Item(4, 6).Value = 1

$ws2.Cells.Item(2, 7).Value = "MSG_UC_MASTER_GET_RUN_UC_TYPE"
$ws2.Cells.Item(3, 7).Value = "MSG_UC_MASTER_GET_RUN_UC_STATE_ID"
$ws2.Cells.Item(4, 7).Value = "MSG_UC_MASTER_GET_RUN_UC_STATE_ACTION"

$ws2.Range("G7").Select()

# --- 3. Populate "ccpd_post_st" --------------------------------------------
$ws = $wb.Worksheets.Item("ccpd_post_st")

for ($i = 0; $i -lt $header.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $header[$i]
}
$ws.Range("A1:G1").Font.Bold = $true

$ws.Cells.Item(2, 1).Value = "usecase/uc_ccpd.c"

$fields = @("begin_time_ms", "mainpump_st", "battery_st", "ht_inlet_st", "ht_outlet_st", `
            "ht_hpad1_st", "ht_hpad2_st", "heater_func_st", "turb_top_st", "turb_side_st", `
            "level_st", "perit_st", "airpressure_st", "eeprom_st", "fan_st", "end_time_ms")
for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $fields[$i]
}

$ws.Cells.Item(2, 2).Value = "CCPD_POST_ST_STRUCT"

$datatypes = @("uint32_t", "uint8_t", "uint8_t", "uint8_t", "uint8_t", `
               "uint8_t", "uint8_t", "uint8_t", "uint8_t", "uint8_t", `
               "uint8_t", "uint8_t", "uint8_t", "uint8_t", "uint8_t", "uint32_t")
for ($i = 0; $i -lt $datatypes.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $datatypes[$i]
}

$bytes = @(4, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 4)
for ($i = 0; $i -lt $bytes.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $bytes[$i]
}
$ws.Cells.Item(2, 3).Value = 110

# NOTE: rows 12/13 (level_st / perit_st) intentionally carry each other's
# "MSG_..." name, matching the source workbook (an upstream naming quirk
# preserved verbatim, not a transcription error).
$msgs = @(
    "MSG_UC_MASTER_GET_CCPD_POST_BEGIN_MS",
    "MSG_UC_MASTER_GET_CCPD_POST_MAINPUMP",
    "MSG_UC_MASTER_GET_CCPD_POST_BATTERY",
    "MSG_UC_MASTER_GET_CCPD_POST_HT_INLET",
    "MSG_UC_MASTER_GET_CCPD_POST_HT_OUTLET",
    "MSG_UC_MASTER_GET_CCPD_POST_HT_HPAD1",
    "MSG_UC_MASTER_GET_CCPD_POST_HT_HPAD2",
    "MSG_UC_MASTER_GET_CCPD_POST_HEATER_HEAT",
    "MSG_UC_MASTER_GET_CCPD_POST_TURB_TOP",
    "MSG_UC_MASTER_GET_CCPD_POST_TURB_SIDE",
    "MSG_UC_MASTER_GET_CCPD_POST_PERIT",
    "MSG_UC_MASTER_GET_CCPD_POST_LEVEL",
    "MSG_UC_MASTER_GET_CCPD_POST_AIRPA",
    "MSG_UC_MASTER_GET_CCPD_POST_EEPROM",
    "MSG_UC_MASTER_GET_CCPD_POST_FAN",
    "MSG_UC_MASTER_GET_CCPD_POST_END_MS"
)
for ($i = 0; $i -lt $msgs.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $msgs[$i]
}

$ws.Range("E17").Select()

# --- 4. Make "ccpd_post_st" the active tab ---------------------------------
$ws.Activate()
$ws.Range("E17").Select()
